# Updating Mental Health Ontology mapping to LSRs.xlsx
# - Rename several "Variable to extract" (column E) codes to more
#   readable/underscored forms (and tweak a couple of column F labels).
# - Append 9 new rows (675-683) describing PTSD3 baseline mappings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column E (and a couple of column F) renames on existing rows
# ---------------------------------------------------------------------

$renames = @(
    @{Row=556; E="Exercise_type";        F="Type of exercise delivered in the intervention arm"}
    @{Row=561; E="Exercise_type";        F="Type of exercise delivered in the intervention arm"}
    @{Row=563; E="Attention_control"}
    @{Row=564; E="intensity_EX"}
    @{Row=565; E="intensity_EX"}
    @{Row=566; E="intervention_length"}
    @{Row=568; E="N_dropout"}
    @{Row=570; E="PTSD3_posttest_mean"}
    @{Row=571; E="PTSD3_posttest_mean"}
    @{Row=572; E="PTSD3_posttest_mean"}
    @{Row=573; E="PTSD3_posttest_n"}
    @{Row=574; E="PTSD3_posttest_n"}
    @{Row=575; E="PTSD3_posttest_n"}
    @{Row=576; E="PTSD3_posttest_sd"}
    @{Row=577; E="PTSD3_posttest_sd"}
    @{Row=578; E="PTSD3_posttest_sd"}
    @{Row=579; E="RoB_PTSD_severity"}
    @{Row=580; E="RoB_PTSD_severity"}
    @{Row=672; E="female_n"}
    @{Row=673; E="PTSD3_questionnaire"}
)

foreach ($item in $renames) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    if ($item.ContainsKey("F")) {
        $ws.Cells.Item($item.Row, 6).Value = $item.F
    }
}

# ---------------------------------------------------------------------
# 2. Append new rows 675-683 (PTSD3 baseline mappings)
# ---------------------------------------------------------------------

$newRows = @(
    @{Row=675; C="2"; E="PTSD3_baseline_n"; F="Number of participants with baseline evaluation on post-traumatic stress symptom"; G="GMHO:0000173"; H="post-traumatic stress symptom severity"; I="Symptom severity relating to a post-traumatic stress symptom."; J="symptom severity"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=676; C="2"; E="PTSD3_baseline_n"; F="Number of participants with baseline evaluation on post-traumatic stress symptom"; G="GMHO:0000209"; H="measurement datum at baseline"; I="Measurement datum that was recorded as baseline data in a study."; J="measurement datum"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=677; C="2"; E="PTSD3_baseline_n"; F="Number of participants with baseline evaluation on post-traumatic stress symptom"; G="GMHO:0000206"; H="number of participants with measurement"; I="Number of intervention participants for whom a measurement was made."; J="number of intervention participants"; K="GMHO:0000173,GMHO:0000209"; L="No Combo"}
    @{Row=678; C="2"; E="PTSD3_baseline_mean"; F="Mean value of baseline post-traumatic stress symptom severity"; G="GMHO:0000173"; H="post-traumatic stress symptom severity"; I="Symptom severity relating to a post-traumatic stress symptom."; J="symptom severity"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=679; C="2"; E="PTSD3_baseline_mean"; F="Mean value of baseline post-traumatic stress symptom severity"; G="GMHO:0000209"; H="measurement datum at baseline"; I="Measurement datum that was recorded as baseline data in a study."; J="measurement datum"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=680; C="2"; E="PTSD3_baseline_mean"; F="Mean value of baseline post-traumatic stress symptom severity"; G="OBI:0000679"; H="average value"; I="A data item that is produced as the output of an averaging data transformation and represents the average value of the input data."; J="data item"; K="GMHO:0000173,GMHO:0000209"; L="No Combo"}
    @{Row=681; C="2"; E="PTSD3_baseline_sd"; F="Standard deviation of baseline post-traumatic stress symptom severity"; G="GMHO:0000173"; H="post-traumatic stress symptom severity"; I="Symptom severity relating to a post-traumatic stress symptom."; J="symptom severity"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=682; C="2"; E="PTSD3_baseline_sd"; F="Standard deviation of baseline post-traumatic stress symptom severity"; G="GMHO:0000209"; H="measurement datum at baseline"; I="Measurement datum that was recorded as baseline data in a study."; J="measurement datum"; K="GMHO:0000087 (Population)"; L="COMBO"; M="GMHO:0000173,GMHO:0000209"}
    @{Row=683; C="2"; E="PTSD3_baseline_sd"; F="Standard deviation of baseline post-traumatic stress symptom severity"; G="OBCS:0000077"; H="standard deviation"; I="A quantitative confidence value that measures the variability of data around the mean."; J="quantitative confidence value"; K="GMHO:0000173,GMHO:0000209"; L="No Combo"}
)

foreach ($item in $newRows) {
    # Column C holds text that looks like a number ("2"); force text via
    # a leading apostrophe so it round-trips as a string, not a number.
    $ws.Cells.Item($item.Row, 3).Value2 = "'" + $item.C
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    $ws.Cells.Item($item.Row, 6).Value = $item.F
    $ws.Cells.Item($item.Row, 7).Value = $item.G
    $ws.Cells.Item($item.Row, 8).Value = $item.H
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
    $ws.Cells.Item($item.Row, 11).Value = $item.K
    $ws.Cells.Item($item.Row, 12).Value = $item.L
    if ($item.ContainsKey("M")) {
        $ws.Cells.Item($item.Row, 13).Value = $item.M
    }
}
